$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83: column B was stored as text "3"; fix it to be a true numeric 3
$ws.Cells.Item(83, 2).Value = 3

# Row 84: new annotation row appended by Ruilin
$ws.Cells.Item(84, 1).Value = "Ruilin"

# B84 must be stored as text "3" (matching source data, not a real number)
$ws.Cells.Item(84, 2).NumberFormat = "@"
$ws.Cells.Item(84, 2).Value = "3"
$ws.Cells.Item(84, 2).ClearFormats()

$ws.Cells.Item(84, 3).Value = "无"
$ws.Cells.Item(84, 4).Value = "FBK"
$ws.Cells.Item(84, 5).Value = "WRI"
$ws.Cells.Item(84, 6).Value = "23ce80a1-f5c9-4d52-8c77-e985ea50fb2a"
$ws.Cells.Item(84, 7).Value = "H1uR4GZRZ_annotated.xlsx"
$ws.Cells.Item(84, 8).Value = "Per your suggestions, we have improved the draft by running additional experiments."
